# Add a new "2021" column (column R) to the table, mirroring the
# formatting of the existing "2020" column (Q), and update the sheet's
# selection to follow the new last column, as Excel would do after a
# user types a new column of data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column header row (row 4): R4 = 2021, same style as Q4 (2020) ---
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("R4").Value = 2021

# --- Data row (row 5): R5 = 42.9, same style as Q5 (47.4) ---
$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("R5").Value = 42.9

$ws.Application.CutCopyMode = $false

# Move the active selection one column to the right, matching the
# author's recorded selection of R9 (was Q9) after editing the table.
[void]$ws.Range("R9").Select()
